$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Matemática Básica"
$ws.Range("C6").Value = "Operações"
$ws.Range("D6").Value = "Tabuada de multiplicação de 6 a 8"
$ws.Range("E6").Value = "<ul>`n`t<li>6 x 7 = 42</li>`n`t<li>6 x 8 = 48</li>`n`t<li>7 x 8 = 56</li>`n</ul>"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
